$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old data went down to row 12)
$ws.Rows("6:12").Delete()

# Header row (Meal / Mon / Thur / Tues / Weds)
$ws.Range("B1").Value = "Meal"
$ws.Range("C1").Value = "Mon"
$ws.Range("D1").Value = "Thur"
$ws.Range("E1").Value = "Tues"
$ws.Range("F1").Value = "Weds"

# Match the bold/bordered header style used by B1/C1 on the new header cells too
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

# "Meal" label column (B2:B5)
$ws.Range("B2").Value = "Breakfast"
$ws.Range("B3").Value = "Lunch"
$ws.Range("B4").Value = "Dinner"
$ws.Range("B5").Value = "Midnight Snack"

# Mon column (C2:C5)
$ws.Range("C2").Value = "Toast"
$ws.Range("C3").Value = "Soup"
$ws.Range("C4").Value = "Curry"
$ws.Range("C5").Value = "Shmores"

# Thur column (D2:D5)
$ws.Range("D2").Value = "Toast"
$ws.Range("D3").Value = "Hotpot"
$ws.Range("D4").Value = "Curry"
$ws.Range("D5").Value = "Chocolate"

# Tues column (E2:E5)
$ws.Range("E2").Value = "Toast"
$ws.Range("E3").Value = "Something Different!"
$ws.Range("E4").Value = "Curry"
$ws.Range("E5").Value = "Shmores"

# Weds column (F2:F5)
$ws.Range("F2").Value = "Toast"
$ws.Range("F3").Value = "Soup"
$ws.Range("F4").Value = "Curry"
$ws.Range("F5").Value = "Biscuits"
